# "Flow Lu Crm Opportunities is added"
#
# Summary of the change being applied:
#  - Duplicate the "LoginTest" sheet into a new "LoginTestFlow" sheet.
#  - Add a brand new "AddOpportunityTestFlow" sheet (right after "test_suite")
#    with a small opportunity-creation dataset (incl. two hyperlinks).
#  - Move "LoginTestFlow" so it sits right after "AddOpportunityTestFlow".
#  - Update the "test_suite" control sheet: drop the old RCRM rows and
#    replace them with the two new "*Flow" test rows.
#  - Relocate the old RCRM control rows (LoginTestReallyCRM,
#    AddAccountTestRCRM, OpportunityTestRCRM) onto "Sheet1", plus a new row
#    for "AddOpportunityTestFlow".
#  - Tidy up a couple of leftover cursor/selection positions on other sheets.
#
# NOTE: worksheet handles returned by Worksheets.Item(...) behave like
# positional handles, not stable object identities - inserting/moving a
# sheet elsewhere in the workbook can silently repoint an old handle at a
# different sheet. So every sheet reference below is re-fetched by name
# immediately before it is used, rather than being cached across any
# Add/Copy/Move call.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate "LoginTest" -> "LoginTestFlow"
# ---------------------------------------------------------------------
$wb.Worksheets.Item("LoginTest").Copy($null, $wb.Worksheets.Item("LoginTest")) | Out-Null
$wb.Worksheets.Item("LoginTest (2)").Name = "LoginTestFlow"

# ---------------------------------------------------------------------
# 2. Add the new "AddOpportunityTestFlow" sheet right after "test_suite"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("test_suite"))
$newSheet.Name = "AddOpportunityTestFlow"

# Move the login-flow copy so it lands right after the new opportunity sheet
$wb.Worksheets.Item("LoginTestFlow").Move($null, $wb.Worksheets.Item("AddOpportunityTestFlow")) | Out-Null

# ---------------------------------------------------------------------
# 3. Populate "AddOpportunityTestFlow"
# ---------------------------------------------------------------------
$oppFlow = $wb.Worksheets.Item("AddOpportunityTestFlow")

$oppFlow.Range("A1").Value = "title"
$oppFlow.Range("B1").Value = "amount"
$oppFlow.Range("C1").Value = "source"
$oppFlow.Range("D1").Value = "email"
$oppFlow.Range("E1").Value = "phone"
$oppFlow.Range("F1").Value = "website"
$oppFlow.Range("G1").Value = "oraganization"
$oppFlow.Range("H1").Value = "contactperson"

$oppFlow.Range("A2").Value = "None"
$oppFlow.Range("B2").Value = 100
$oppFlow.Range("C2").Value = "Partner"
$oppFlow.Range("D2").Value = "a@gmail.com"
$oppFlow.Hyperlinks.Add($oppFlow.Range("D2"), "mailto:a@gmail.com") | Out-Null
$oppFlow.Range("D2").Style = "Hyperlink"
$oppFlow.Range("E2").Value = 900000
$oppFlow.Range("F2").Value = "www.google.com"
$oppFlow.Hyperlinks.Add($oppFlow.Range("F2"), "http://www.google.com/") | Out-Null
$oppFlow.Range("F2").Style = "Hyperlink"
$oppFlow.Range("G2").Value = "iskcon"
$oppFlow.Range("H2").Value = "bluetoothprabu"

$oppFlow.PageSetup.Orientation = 1
$oppFlow.Range("H1").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Rework "test_suite": drop old rows, add the two new *Flow rows
# ---------------------------------------------------------------------
$testSuite = $wb.Worksheets.Item("test_suite")
$testSuite.Range("A2:B4").ClearContents() | Out-Null
$testSuite.Range("A2").Value = "LoginTestFlow"
$testSuite.Range("B2").Value = "Y"
$testSuite.Range("A3").Value = "AddOpportunityTestFlow"
$testSuite.Range("B3").Value = "Y"
$testSuite.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. Relocate the old RCRM control rows onto "Sheet1"
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("A10").Value = "LoginTestReallyCRM"
$sheet1.Range("B10").Value = "Y"
$sheet1.Range("A11").Value = "AddAccountTestRCRM"
$sheet1.Range("B11").Value = "Y"
$sheet1.Range("A12").Value = "OpportunityTestRCRM"
$sheet1.Range("B12").Value = "Y"
$sheet1.Range("A13").Value = "AddOpportunityTestFlow"
$sheet1.Range("B13").Value = "Y"
$sheet1.Range("H15").Select() | Out-Null

# ---------------------------------------------------------------------
# 6. Leftover cursor/selection tidy-up
# ---------------------------------------------------------------------
$wb.Worksheets.Item("LoginTest").Range("A1:C2").Select() | Out-Null

$wb.Worksheets.Item("AddOpportunityTestFlow").Activate() | Out-Null
